$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q4" sheet right before the existing "2022-Q2" sheet.
#    Duplicate "2022-Q2" (same column layout/header/styles) and overwrite the
#    values in place so the original cell formatting (borders, bold header,
#    centered index column) carries over untouched.
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2022-Q2")
$refSheet.Copy($refSheet)
$q4 = $wb.Worksheets.Item("2022-Q2 (2)")
$q4.Name = "2022-Q4"

# Header row text is unchanged (same columns as every other quarter sheet) -
# only the data rows need to be replaced.

# Data rows (A stays numeric index, B is text fund code, C is text name,
# D/E/F/G are text-preserving numeric-looking strings, H is a real number)
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'005457"
$q4.Range("C2").Value = "景顺长城量化小盘股票"
$q4.Range("D2").Value = "'5.08"
$q4.Range("E2").Value = "'94.36"
$q4.Range("F2").Value = "'1.46"
$q4.Range("G2").Value = "'0.0742"
$q4.Range("H2").Value = 5

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'008851"
$q4.Range("C3").Value = "景顺长城量化对冲策略三个月定期开放灵活配置混合"
$q4.Range("D3").Value = "'2.37"
$q4.Range("E3").Value = "'71.22"
$q4.Range("F3").Value = "'1.11"
$q4.Range("G3").Value = "'0.0263"
$q4.Range("H3").Value = 5

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "'015496"
$q4.Range("C4").Value = "景顺中证1000指数增强C"
$q4.Range("D4").Value = "'0.86"
$q4.Range("E4").Value = "'92.30"
$q4.Range("F4").Value = "'1.49"
$q4.Range("G4").Value = "'0.0128"
$q4.Range("H4").Value = 5

# The reference sheet only had 4 data rows - this is the 5th, so its "A"
# cell (the bordered/centered index-column style) needs to be stamped onto
# the new row before the value is written over it.
$q4.Range("A2").Copy()
$q4.Range("A5").PasteSpecial(-4122)
$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "'015495"
$q4.Range("C5").Value = "景顺中证1000指数增强A"
$q4.Range("D5").Value = "'0.67"
$q4.Range("E5").Value = "'92.30"
$q4.Range("F5").Value = "'1.49"
$q4.Range("G5").Value = "'0.0100"
$q4.Range("H5").Value = 5

# ---------------------------------------------------------------------------
# 2) "总计" sheet: insert a new row 2 for the 2022-Q4 summary and shift the
#    rest of the quarters down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# The inserted row inherits formatting from the row above (the header) -
# strip that off the data cells, then re-apply the plain index-column style
# that every other data row in column A uses.
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.12

# Renumber the index column (A) sequentially 0..7 for the 8 data rows now
# that the new quarter has been inserted at the top.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
$total.Range("A9").Value = 7
